$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column S (year 2023) values, keyed by row number.
$values = @{
    3  = 2023
    4  = 1926.4
    5  = 1929.2
    6  = 24982
    7  = 24520
    8  = 1481.1
    9  = 1068.5
    10 = 443
    11 = 860.8
    12 = 240.1
    13 = 1057.7
    14 = 1
}

foreach ($r in 3..14) {
    # Copy the formatting from the existing last column (R) of the row onto
    # the new S cell before writing its value, so the new column matches the
    # look of the rest of the table.
    $ws.Range("R$r").Copy()
    $ws.Range("S$r").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    $ws.Range("S$r").Value = $values[$r]
}

# Move the active selection to match the edited workbook.
[void]$ws.Range("J19").Select()
